$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new departure rows (21, 22) for recurring LOT flights on Monday, Jan 09
# Row 21: flight LO3802 (2:55 PM scheduled)
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Monday, Jan 09"
$ws.Cells.Item(21, 3).Value = "2:55 PM"
$ws.Cells.Item(21, 4).Value = "LO3802"
$ws.Cells.Item(21, 5).Value = "Warsaw"
$ws.Cells.Item(21, 6).Value = "(WAW)"
$ws.Cells.Item(21, 7).Value = "LOT "
$ws.Cells.Item(21, 8).Value = "E75S"
$ws.Cells.Item(21, 9).Value = "(SP-LIA)"
$ws.Cells.Item(21, 10).Value = "2:58 PM"
$ws.Cells.Item(21, 12).Value = "0 hours, 3 minutes"

# Row 22: flight LO3808 (4:45 PM scheduled)
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Monday, Jan 09"
$ws.Cells.Item(22, 3).Value = "4:45 PM"
$ws.Cells.Item(22, 4).Value = "LO3808"
$ws.Cells.Item(22, 5).Value = "Warsaw"
$ws.Cells.Item(22, 6).Value = "(WAW)"
$ws.Cells.Item(22, 7).Value = "LOT "
$ws.Cells.Item(22, 8).Value = "E75S"
$ws.Cells.Item(22, 9).Value = "(SP-LID)"
$ws.Cells.Item(22, 10).Value = "5:12 PM"
$ws.Cells.Item(22, 12).Value = "0 hours, 27 minutes"
